$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Delete the row containing "Implementing game-renderer" (row 4)
$ws.Rows.Item(4).Delete()

# After the first deletion, "Implementing one static game field" has shifted
# up from row 13 to row 12 - delete it too.
$ws.Rows.Item(12).Delete()
